$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of date range) ---
$ws.Range("A8").Characters(21, 1).Text = "8"
$ws.Range("C9").Characters(27, 9).Text = "2/19/2024"
$ws.Range("C9").Characters(47, 9).Text = "2/25/2024"

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("J15").Value = 8
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -45.454545454545
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 39
$ws.Range("K16").Value = -15.384615384615
$ws.Range("L16").Value = 43.478260869565
$ws.Range("M16").Value = -8.333333333333
$ws.Range("N16").Value = -83.076923076923
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("I17").Value = 34
$ws.Range("J17").Value = 28
$ws.Range("K17").Value = 21.428571428571
$ws.Range("L17").Value = 17.241379310344
$ws.Range("M17").Value = 126.666666666667
$ws.Range("N17").Value = -10.526315789473
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 12.5
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = 12.5
$ws.Range("L18").Value = 16.129032258064
$ws.Range("M18").Value = -5.263157894736
$ws.Range("N18").Value = -84.745762711864
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 108
$ws.Range("J19").Value = 107
$ws.Range("K19").Value = 0.934579439252
$ws.Range("L19").Value = -3.571428571428
$ws.Range("M19").Value = 77.049180327868
$ws.Range("N19").Value = -12.195121951219
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -23.076923076923
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -35.714285714285
$ws.Range("L20").Value = -25
$ws.Range("M20").Value = -51.351351351351
$ws.Range("N20").Value = -94.321766561514
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -37.931034482758
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = -16
$ws.Range("I21").Value = 229
$ws.Range("J21").Value = 242
$ws.Range("K21").Value = -5.371900826446
$ws.Range("L21").Value = 3.153153153153
$ws.Range("M21").Value = 21.164021164021
$ws.Range("N21").Value = -74.890350877193
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -57.142857142857
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = -45.454545454545
$ws.Range("M22").Value = -25
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 25.806451612903
$ws.Range("F24").Value = 178
$ws.Range("G24").Value = 138
$ws.Range("H24").Value = 28.985507246376
$ws.Range("I24").Value = 325
$ws.Range("J24").Value = 238
$ws.Range("K24").Value = 36.554621848739
$ws.Range("L24").Value = 96.969696969697
$ws.Range("M24").Value = 218.627450980392
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 85.714285714285
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 51.515151515151
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 88
$ws.Range("K25").Value = -5.681818181818
$ws.Range("L25").Value = -5.681818181818
$ws.Range("M25").Value = 6.410256410256
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = -80
$ws.Range("L26").Value = -33.333333333333
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -54.545454545454
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -41.176470588235
$ws.Range("L27").Value = 11.111111111111

# --- Cells changing from text placeholder ("0"/"***.*" ) to numeric ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 1
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("M28").Value = -100
$ws.Range("M28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M29").Value = -100
$ws.Range("M29").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Cells changing from numeric to text placeholder ("0"/"***.*") ---
# Use the fully-text row 23 as a format donor so the resulting style
# matches the workbook's existing "General/text" style exactly.
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D23").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E23").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D23").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E23").Copy()
$ws.Range("E30").PasteSpecial(-4122)